# Applies the "can indirectly update EntryPoints and PointDefs" change.
# This re-runs a simulated sync/update cycle against the fixture data:
#   - several rows across Defs / Point Defs / Entry Base / Entry Points get
#     touched (new _uid + refreshed _created/_updated timestamps),
#   - some rows get relabelled / retyped / toggled,
#   - two brand-new rows appear (one in "Point Defs", one in "Entry Points").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Defs
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Defs")

$ws.Range("A2").Value = "lgs5e3pj-elk6"
$ws.Range("B2").Value = "2023-04-22T15:41:12.391Z"
$ws.Range("C2").Value = "lgs5e3pj"
$ws.Range("D2").Value = "'FALSE"

$ws.Range("A3").Value = "lgs5e3pj-6tewf"
$ws.Range("B3").Value = "2023-04-22T15:41:12.391Z"
$ws.Range("C3").Value = "lgs5e3q3"
$ws.Range("D3").Value = "'TRUE"

$ws.Range("A4").Value = "lgs5e3pj-1ve7"
$ws.Range("B4").Value = "2023-04-22T15:41:12.391Z"
$ws.Range("C4").Value = "lgs5e3pj"

$ws.Range("A5").Value = "lgs5e3q3-04qv"
$ws.Range("B5").Value = "2023-04-22T15:41:12.391Z"
$ws.Range("C5").Value = "lgs5e3q3"
$ws.Range("E5").Value = "ay7l"
$ws.Range("F5").Value = "Two Relabeled"
$ws.Range("G5").Value = "2️⃣"
$ws.Range("H5").Value = "Scoped at an **hour**, cause why not have that option?"
$ws.Range("I5").Value = "HOUR"

# ---------------------------------------------------------------------------
# Sheet: Point Defs
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Point Defs")

$ws.Range("A2").Value = "lgs5e3pk-0cjl"
$ws.Range("B2").Value = "2023-04-22T15:41:12.392Z"
$ws.Range("C2").Value = "lgs5e3pk"
$ws.Range("F2").Value = "0eze"

$ws.Range("A3").Value = "lgs5e3pk-d89q"
$ws.Range("B3").Value = "2023-04-22T15:41:12.392Z"
$ws.Range("C3").Value = "lgs5e3q3"
$ws.Range("D3").Value = $true

$ws.Range("A4").Value = "lgs5e3pk-0kt2"
$ws.Range("B4").Value = "2023-04-22T15:41:12.392Z"
$ws.Range("C4").Value = "lgs5e3q3"

$ws.Range("A5").Value = "lgs5e3q3-3s9h"
$ws.Range("B5").Value = "2023-04-22T15:41:12.392Z"
$ws.Range("C5").Value = "lgs5e3q3"
$ws.Range("F5").Value = "0pc6"
$ws.Range("G5").Value = "Test Relabel"
$ws.Range("H5").Value = "#️⃣"
$ws.Range("I5").Value = "Set a description"
$ws.Range("J5").Value = "NUM"

# New row: a second point def now points back at "Boolean Thing" / 0tb7
$ws.Range("A6").Value = "lgs5e3q3-xdhk"
$ws.Range("B6").Value = "2023-04-22T15:41:12.392Z"
$ws.Range("C6").Value = "lgs5e3q3"
$ws.Range("D6").Value = $false
$ws.Range("E6").Value = "ay7l"
$ws.Range("F6").Value = "0tb7"
$ws.Range("G6").Value = "Boolean Thing"
$ws.Range("H6").Value = "👎"
$ws.Range("I6").Value = "Orig desc"
$ws.Range("J6").Value = "BOOL"
$ws.Range("K6").Value = "COUNT"

# ---------------------------------------------------------------------------
# Sheet: Entry Base
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Entry Base")

$ws.Range("A2").Value = "lgs5e3pk-ust9"
$ws.Range("B2").Value = "2023-04-22T15:41:12.392Z"
$ws.Range("C2").Value = "lgs5e3q3"

$ws.Range("A3").Value = "lgs5e3pk-euus"
$ws.Range("B3").Value = "2023-04-22T15:41:12.392Z"
$ws.Range("C3").Value = "lgs5e3pk"
$ws.Range("F3").Value = "lgs5e3pv-5ph5n"
$ws.Range("G3").Value = "2023-04-22T10:41:12"

$ws.Range("A4").Value = "lgs5e3q3-0jn9"
$ws.Range("B4").Value = "2023-04-22T15:41:12.392Z"
$ws.Range("C4").Value = "lgs5e3q3"

# ---------------------------------------------------------------------------
# Sheet: Entry Points
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Entry Points")

$ws.Range("A2").Value = "lgs5e3pk-00bg"
$ws.Range("B2").Value = "2023-04-22T15:41:12.392Z"
$ws.Range("C2").Value = "lgs5e3q3"
$ws.Range("D2").Value = $true

$ws.Range("A3").Value = "lgs5e3pk-5gq2"
$ws.Range("B3").Value = "2023-04-22T15:41:12.392Z"
$ws.Range("C3").Value = "lgs5e3q3"

$ws.Range("A4").Value = "lgs5e3q3-0y6e"
$ws.Range("B4").Value = "2023-04-22T15:41:12.392Z"
$ws.Range("C4").Value = "lgs5e3q3"
$ws.Range("F4").Value = "0tb7"
$ws.Range("H4").Value = "'true"

# New row: the numeric point (0pc6) now also has an entry-point value of 6
$ws.Range("A5").Value = "lgs5e3q3-065z"
$ws.Range("B5").Value = "2023-04-22T15:41:12.392Z"
$ws.Range("C5").Value = "lgs5e3q3"
$ws.Range("D5").Value = $false
$ws.Range("E5").Value = "ay7l"
$ws.Range("F5").Value = "0pc6"
$ws.Range("G5").Value = "lgricx7k-08al"
$ws.Range("H5").Value = "'6"
